$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.451.16"
$ws.Range("E2").Value = "  +0.49%  "

$ws.Range("D3").Value = "3.335.56"
$ws.Range("E3").Value = "  +0.16%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").Value = "190.18"
$ws.Range("E5").Value = "  +4.70%  "

$ws.Range("D6").Value = "566.13"
$ws.Range("E6").Value = "  +1.84%  "

$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.07%  "

$ws.Range("D8").Value = "0.590"
$ws.Range("E8").Value = "  -0.02%  "

$ws.Range("D9").Value = "3.327.59"
$ws.Range("E9").Value = "  +0.16%  "

$ws.Range("D10").Value = "0.186"
$ws.Range("E10").Value = "  +1.84%  "

$ws.Range("D11").Value = "0.592"
$ws.Range("E11").Value = "  +1.37%  "

$ws.Range("E12").Value = "  +1.67%  "

$ws.Range("D13").Value = "0.0000274"
$ws.Range("E13").Value = "  +3.96%  "

$ws.Range("D14").Value = "8.72"
$ws.Range("E14").Value = "  +1.40%  "

$ws.Range("D15").Value = "3.880.80"
$ws.Range("E15").Value = "  +0.66%  "

$ws.Range("D16").Value = "611.83"
$ws.Range("E16").Value = "  +1.82%  "

$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "66.550.96"
$ws.Range("E17").Value = "  +0.63%  "

$ws.Range("B18").Value = "Chainlink"
$ws.Range("C18").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D18").Value = "18.15"
$ws.Range("E18").Value = "  +1.38%  "

$ws.Range("E19").Value = "  +1.28%  "

$ws.Range("D20").Value = "3.343.77"
$ws.Range("E20").Value = "  +0.18%  "

$ws.Range("D21").Value = "11.24"
$ws.Range("E21").Value = "  -1.09%  "

$ws.Range("D22").Value = "0.917"
$ws.Range("E22").Value = "  +1.75%  "

$ws.Range("D23").Value = "18.68"
$ws.Range("E23").Value = "  +11.73%  "

$ws.Range("D24").Value = "5.19"
$ws.Range("E24").Value = "  +2.61%  "

$ws.Range("D25").Value = "101.69"
$ws.Range("E25").Value = "  +2.90%  "

$ws.Range("D26").Value = "4.05"
$ws.Range("E26").Value = "  +1.00%  "

$ws.Range("D27").Value = "2.78"
$ws.Range("E27").Value = "  +3.55%  "

$ws.Range("D28").Value = "9.80"
$ws.Range("E28").Value = "  +5.49%  "

$ws.Range("D29").Value = "8.74"
$ws.Range("E29").Value = "  +0.97%  "

$ws.Range("D30").Value = "30.57"
$ws.Range("E30").Value = "  +1.05%  "

$ws.Range("E31").Value = "  +9.44%  "

$ws.Range("D32").Value = "4.05"
$ws.Range("E32").Value = "  +8.08%  "

$ws.Range("D33").Value = "567.61"
$ws.Range("E33").Value = "  +3.30%  "

$ws.Range("E34").Value = "  +1.60%  "

$ws.Range("E35").Value = "  +2.01%  "

$ws.Range("D36").Value = "3.733.92"
$ws.Range("E36").Value = "  -2.44%  "

$ws.Range("D37").Value = "57.41"
$ws.Range("E37").Value = "  -0.30%  "

$ws.Range("D38").Value = "0.999"
$ws.Range("E38").Value = "  -0.02%  "

$ws.Range("D39").Value = "0.0₃0733"
$ws.Range("E39").Value = "  +4.99%  "

$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").Value = "0.132"
$ws.Range("E40").Value = "  +6.08%  "

$ws.Range("B41").Value = "InjectiveProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D41").Value = "34.18"
$ws.Range("E41").Value = "  +7.14%  "

$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").Value = "3.32"
$ws.Range("E42").Value = "  -1.77%  "

$ws.Range("D43").Value = "2.72"
$ws.Range("E43").Value = "  +3.10%  "

$ws.Range("D44").Value = "3.41"
$ws.Range("E44").Value = "  +10.10%  "

$ws.Range("E45").Value = "  +1.04%  "

$ws.Range("D46").Value = "0.0429"
$ws.Range("E46").Value = "  +4.45%  "

$ws.Range("D47").Value = "3.25"
$ws.Range("E47").Value = "  +4.78%  "

$ws.Range("E48").Value = "  +0.85%  "

$ws.Range("D49").Value = "2.62"
$ws.Range("E49").Value = "  +0.25%  "

$ws.Range("E50").Value = "  +0.25%  "

$ws.Range("E51").Value = "  +4.18%  "
